$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four "Resolving-Mac" sending-cluster rows (old rows 10-13)
$ws.Rows.Item(10).Resize(4).EntireRow.Delete()

# Row 2: updated TPM-derived values
$ws.Range("I2").Value = 0.5062291280850276
$ws.Range("J2").Value = 0.5062291280850276
$ws.Range("M2").Value = 10.13228833333333
$ws.Range("N2").Value = 30.396865
$ws.Range("O2").Value = 0.4620357097718242
$ws.Range("P2").Value = 0.4620357097718243
$ws.Range("Q2").Value = 2.937728659931111
$ws.Range("R2").Value = 26.43955793938
$ws.Range("S2").Value = 0.2338959345019375
$ws.Range("T2").Value = 0.2338959345019375

# Row 3: updated TPM-derived values
$ws.Range("I3").Value = 0.5062291280850276
$ws.Range("J3").Value = 0.5062291280850276
$ws.Range("O3").Value = 0.3432940418074481
$ws.Range("P3").Value = 0.3432940418074482
$ws.Range("S3").Value = 0.1737854434609695
$ws.Range("T3").Value = 0.1737854434609695

# Row 4: updated TPM-derived values
$ws.Range("I4").Value = 0.5062291280850276
$ws.Range("J4").Value = 0.5062291280850276
$ws.Range("M4").Value = 1.155416666666667
$ws.Range("N4").Value = 3.46625
$ws.Range("O4").Value = 0.05268738335340128
$ws.Range("P4").Value = 0.05268738335340128
$ws.Range("Q4").Value = 0.3349984272222222
$ws.Range("R4").Value = 3.014985845
$ws.Range("S4").Value = 0.02667188813607393
$ws.Range("T4").Value = 0.02667188813607393

# Row 5: updated TPM-derived values
$ws.Range("I5").Value = 0.5062291280850276
$ws.Range("J5").Value = 0.5062291280850276
$ws.Range("M5").Value = 3.113636666666667
$ws.Range("N5").Value = 9.340910000000001
$ws.Range("O5").Value = 0.1419828650673262
$ws.Range("P5").Value = 0.1419828650673262
$ws.Range("Q5").Value = 0.9027595121022223
$ws.Range("R5").Value = 8.124835608920002
$ws.Range("S5").Value = 0.07187586198604667
$ws.Range("T5").Value = 0.07187586198604669

# Row 6: updated TPM-derived values
$ws.Range("G6").Value = 0.282802
$ws.Range("H6").Value = 0.848406
$ws.Range("I6").Value = 0.4937708719149724
$ws.Range("J6").Value = 0.4937708719149724
$ws.Range("M6").Value = 10.13228833333333
$ws.Range("N6").Value = 30.396865
$ws.Range("O6").Value = 0.4620357097718242
$ws.Range("P6").Value = 0.4620357097718243
$ws.Range("Q6").Value = 2.865431405243333
$ws.Range("R6").Value = 25.78888264719
$ws.Range("S6").Value = 0.2281397752698868
$ws.Range("T6").Value = 0.2281397752698868

# Row 7: updated TPM-derived values
$ws.Range("G7").Value = 0.282802
$ws.Range("H7").Value = 0.848406
$ws.Range("I7").Value = 0.4937708719149724
$ws.Range("J7").Value = 0.4937708719149724
$ws.Range("O7").Value = 0.3432940418074481
$ws.Range("P7").Value = 0.3432940418074482
$ws.Range("Q7").Value = 2.129024895313333
$ws.Range("R7").Value = 19.16122405782
$ws.Range("S7").Value = 0.1695085983464787
$ws.Range("T7").Value = 0.1695085983464787

# Row 8: updated TPM-derived values
$ws.Range("G8").Value = 0.282802
$ws.Range("H8").Value = 0.848406
$ws.Range("I8").Value = 0.4937708719149724
$ws.Range("J8").Value = 0.4937708719149724
$ws.Range("M8").Value = 1.155416666666667
$ws.Range("N8").Value = 3.46625
$ws.Range("O8").Value = 0.05268738335340128
$ws.Range("P8").Value = 0.05268738335340128
$ws.Range("Q8").Value = 0.3267541441666667
$ws.Range("R8").Value = 2.9407872975
$ws.Range("S8").Value = 0.02601549521732735
$ws.Range("T8").Value = 0.02601549521732735

# Row 9: updated TPM-derived values
$ws.Range("G9").Value = 0.282802
$ws.Range("H9").Value = 0.848406
$ws.Range("I9").Value = 0.4937708719149724
$ws.Range("J9").Value = 0.4937708719149724
$ws.Range("M9").Value = 3.113636666666667
$ws.Range("N9").Value = 9.340910000000001
$ws.Range("O9").Value = 0.1419828650673262
$ws.Range("P9").Value = 0.1419828650673262
$ws.Range("Q9").Value = 0.8805426766066667
$ws.Range("R9").Value = 7.924884089460001
$ws.Range("S9").Value = 0.07010700308127954
$ws.Range("T9").Value = 0.07010700308127955
